# Scheduled runner update: refresh market-price derived columns (H-N) per Leve row
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 995.375
$ws.Range("I43").Value = 1036.2
$ws.Range("J43").Value = 927.3333
$ws.Range("K43").Value = 1036.2
$ws.Range("L43").Value = 927.3333
$ws.Range("M43").Value = -967.2
$ws.Range("N43").Value = -1065.3333
$ws.Range("H51").Value = 2082.923
$ws.Range("I51").Value = 2065.4443
$ws.Range("J51").Value = 2122.25
$ws.Range("K51").Value = 2065.4443
$ws.Range("L51").Value = 2122.25
$ws.Range("M51").Value = -1581.4443
$ws.Range("N51").Value = -3090.25
$ws.Range("H134").Value = 51000
$ws.Range("J134").Value = 51000
$ws.Range("L134").Value = 51000
$ws.Range("N134").Value = -61140
$ws.Range("H137").Value = 1286
$ws.Range("I137").Value = 953.875
$ws.Range("J137").Value = 1665.5714
$ws.Range("K137").Value = 2861.625
$ws.Range("L137").Value = 4996.7142
$ws.Range("M137").Value = -311.625
$ws.Range("N137").Value = -10096.7142
$ws.Range("H138").Value = 10419251
$ws.Range("I138").Value = 2843.6316
$ws.Range("J138").Value = 17243794
$ws.Range("K138").Value = 8530.8948
$ws.Range("L138").Value = 51731382
$ws.Range("M138").Value = -3390.8948
$ws.Range("N138").Value = -51741662

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2500.5178
$ws.Range("I32").Value = 1716.7273
$ws.Range("J32").Value = 5374.4165
$ws.Range("K32").Value = 1716.7273
$ws.Range("L32").Value = 5374.4165
$ws.Range("M32").Value = -1429.7273
$ws.Range("N32").Value = -5948.4165
$ws.Range("H61").Value = 1952.6177
$ws.Range("I61").Value = 1138.32
$ws.Range("J61").Value = 4214.5557
$ws.Range("K61").Value = 1138.32
$ws.Range("L61").Value = 4214.5557
$ws.Range("M61").Value = -926.3199999999999
$ws.Range("N61").Value = -4638.5557
$ws.Range("H132").Value = 3065.32
$ws.Range("I132").Value = 2648.5
$ws.Range("J132").Value = 3595.818
$ws.Range("K132").Value = 7945.5
$ws.Range("L132").Value = 10787.454
$ws.Range("M132").Value = -5415.5
$ws.Range("N132").Value = -15847.454
$ws.Range("H136").Value = 1952.6177
$ws.Range("I136").Value = 1138.32
$ws.Range("J136").Value = 4214.5557
$ws.Range("K136").Value = 3414.96
$ws.Range("L136").Value = 12643.6671
$ws.Range("M136").Value = -864.96
$ws.Range("N136").Value = -17743.6671

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 499.91666
$ws.Range("I64").Value = 432.8
$ws.Range("J64").Value = 547.8570999999999
$ws.Range("K64").Value = 432.8
$ws.Range("L64").Value = 547.8570999999999
$ws.Range("M64").Value = -207.8
$ws.Range("N64").Value = -997.8570999999999
$ws.Range("H67").Value = 499.91666
$ws.Range("I67").Value = 432.8
$ws.Range("J67").Value = 547.8570999999999
$ws.Range("K67").Value = 432.8
$ws.Range("L67").Value = 547.8570999999999
$ws.Range("M67").Value = 347.2
$ws.Range("N67").Value = -2107.8571
$ws.Range("H132").Value = 12000
$ws.Range("J132").Value = 12000
$ws.Range("L132").Value = 12000
$ws.Range("N132").Value = -22120
$ws.Range("H134").Value = 2710.65
$ws.Range("I134").Value = 1826.2667
$ws.Range("J134").Value = 5363.8
$ws.Range("K134").Value = 5478.800099999999
$ws.Range("L134").Value = 16091.4
$ws.Range("M134").Value = -2943.800099999999
$ws.Range("N134").Value = -21161.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2219.9697
$ws.Range("I31").Value = 1607.6666
$ws.Range("J31").Value = 2569.8572
$ws.Range("K31").Value = 1607.6666
$ws.Range("L31").Value = 2569.8572
$ws.Range("M31").Value = -1312.6666
$ws.Range("N31").Value = -3159.8572
$ws.Range("H34").Value = 2219.9697
$ws.Range("I34").Value = 1607.6666
$ws.Range("J34").Value = 2569.8572
$ws.Range("K34").Value = 1607.6666
$ws.Range("L34").Value = 2569.8572
$ws.Range("M34").Value = -1405.6666
$ws.Range("N34").Value = -2973.8572

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 898.05884
$ws.Range("I23").Value = 1512
$ws.Range("J23").Value = 352.33334
$ws.Range("K23").Value = 4536
$ws.Range("L23").Value = 1057.00002
$ws.Range("M23").Value = -4301
$ws.Range("N23").Value = -1527.00002
$ws.Range("H68").Value = 1555.3334
$ws.Range("I68").Value = 1435.1082
$ws.Range("J68").Value = 1817
$ws.Range("K68").Value = 4305.3246
$ws.Range("L68").Value = 5451
$ws.Range("M68").Value = -3494.3246
$ws.Range("N68").Value = -7073
$ws.Range("H71").Value = 1555.3334
$ws.Range("I71").Value = 1435.1082
$ws.Range("J71").Value = 1817
$ws.Range("K71").Value = 12915.9738
$ws.Range("L71").Value = 16353
$ws.Range("M71").Value = -8859.9738
$ws.Range("N71").Value = -24465
$ws.Range("H122").Value = 655.9524
$ws.Range("I122").Value = 557.5
$ws.Range("J122").Value = 787.2222
$ws.Range("K122").Value = 5017.5
$ws.Range("L122").Value = 7084.999800000001
$ws.Range("M122").Value = -2567.5
$ws.Range("N122").Value = -11984.9998
$ws.Range("H131").Value = 2496.169
$ws.Range("I131").Value = 388.8889
$ws.Range("K131").Value = 1166.6667
$ws.Range("M131").Value = 3873.3333
$ws.Range("H132").Value = 1048.1
$ws.Range("I132").Value = 564.3333
$ws.Range("J132").Value = 1255.4286
$ws.Range("K132").Value = 5078.9997
$ws.Range("L132").Value = 11298.8574
$ws.Range("M132").Value = -2548.9997
$ws.Range("N132").Value = -16358.8574
$ws.Range("H138").Value = 2915.9167
$ws.Range("I138").Value = 832.3333
$ws.Range("J138").Value = 9166.666999999999
$ws.Range("K138").Value = 2496.9999
$ws.Range("L138").Value = 27500.001
$ws.Range("M138").Value = 2643.0001
$ws.Range("N138").Value = -37780.001
$ws.Range("H140").Value = 2236.875
$ws.Range("I140").Value = 668.46155
$ws.Range("K140").Value = 2005.38465
$ws.Range("M140").Value = 3174.61535
$ws.Range("H141").Value = 4851.1113
$ws.Range("I141").Value = 6132
$ws.Range("J141").Value = 3250
$ws.Range("K141").Value = 18396
$ws.Range("L141").Value = 9750
$ws.Range("M141").Value = -13216
$ws.Range("N141").Value = -20110

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H92").Value = 20000
$ws.Range("J92").Value = 20000
$ws.Range("L92").Value = 20000
$ws.Range("N92").Value = -23744
$ws.Range("H94").Value = 39000
$ws.Range("J94").Value = 39000
$ws.Range("L94").Value = 39000
$ws.Range("N94").Value = -40352

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3051.4827
$ws.Range("I7").Value = 2177
$ws.Range("J7").Value = 3445
$ws.Range("K7").Value = 2177
$ws.Range("L7").Value = 3445
$ws.Range("M7").Value = -2065
$ws.Range("N7").Value = -3669
$ws.Range("H45").Value = 15676.667
$ws.Range("J45").Value = 18494.5
$ws.Range("L45").Value = 18494.5
$ws.Range("N45").Value = -19308.5
$ws.Range("H126").Value = 3051.4827
$ws.Range("I126").Value = 2177
$ws.Range("J126").Value = 3445
$ws.Range("K126").Value = 6531
$ws.Range("L126").Value = 10335
$ws.Range("M126").Value = -4061
$ws.Range("N126").Value = -15275
$ws.Range("H132").Value = 4756
$ws.Range("I132").Value = 4456.3335
$ws.Range("J132").Value = 5295.4
$ws.Range("K132").Value = 13369.0005
$ws.Range("L132").Value = 15886.2
$ws.Range("M132").Value = -10839.0005
$ws.Range("N132").Value = -20946.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 38335
$ws.Range("J4").Value = 55001.5
$ws.Range("L4").Value = 55001.5
$ws.Range("N4").Value = -55227.5
$ws.Range("H124").Value = 39966.668
$ws.Range("J124").Value = 39966.668
$ws.Range("L124").Value = 39966.668
$ws.Range("N124").Value = -49786.668
$ws.Range("H126").Value = 77349.62
$ws.Range("I126").Value = 91340.45
$ws.Range("J126").Value = 400
$ws.Range("K126").Value = 274021.35
$ws.Range("L126").Value = 1200
$ws.Range("M126").Value = -271551.35
$ws.Range("N126").Value = -6140
$ws.Range("H132").Value = 14707228
$ws.Range("I132").Value = 17858024
$ws.Range("J132").Value = 3509
$ws.Range("K132").Value = 53574072
$ws.Range("L132").Value = 10527
$ws.Range("M132").Value = -53571542
$ws.Range("N132").Value = -15587
$ws.Range("H135").Value = 89811.84
$ws.Range("J135").Value = 89811.84
$ws.Range("L135").Value = 89811.84
$ws.Range("N135").Value = -99951.84
